$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "available" column (H) previously held the text "Yes" for every
# banquet row. Replace it with the numeric value 1, formatted as a plain
# (General) number, for each data row (2 through 51).
for ($r = 2; $r -le 51; $r++) {
    $ws.Cells.Item($r, 8).Value = 1
    $ws.Cells.Item($r, 8).NumberFormat = "General"
}
